$wb = $excel.ActiveWorkbook

# Update Ryan's sheet with his logged hours
$ryan = $wb.Worksheets.Item("Ryan")

$ryan.Range("B2").Value = 3
$ryan.Range("C2").Value = 0.5

$ryan.Range("B3").Value = 1
$ryan.Range("C3").Value = 2

$ryan.Range("D4").Value = 1.5

$ryan.Range("D5").Value = 4

$ryan.Range("D6").Value = 3

$ryan.Range("B7").Value = 0.5
$ryan.Range("E7").Value = 1

$ryan.Range("E8").Value = 11

$ryan.Range("E9").Value = 5
$ryan.Range("F9").Value = 2

$ryan.Range("E10").Value = 3
$ryan.Range("F10").Value = 2.5

$ryan.Range("F11").Value = 8

$ryan.Range("F12").Value = 4.5
$ryan.Range("G12").Value = 0.5

$ryan.Range("E13").Value = 2.5
$ryan.Range("F13").Value = 3
$ryan.Range("G13").Value = 2.5

# Selection moves to G10 on Ryan's sheet and it becomes the active/selected tab
$ryan.Range("G10").Select()
$ryan.Activate()

$wb.Save()
